$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price column (D) updates - force text to preserve exact formatting
# (some values look numeric and would otherwise be auto-converted,
# losing trailing zeros / exact decimal representation)
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '43.189.11'
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.282.76'
$ws.Range("D3").Style = "Normal"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '113.43'
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '265.68'
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.622'
$ws.Range("D7").Style = "Normal"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.611'
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '47.78'
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0931'
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '9.12'
$ws.Range("D12").Style = "Normal"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '2.624.47'
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.868'
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.282.43'
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '43.228.06'
$ws.Range("D18").Style = "Normal"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.86'
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '71.59'
$ws.Range("D21").Style = "Normal"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '232.94'
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '9.68'
$ws.Range("D24").Style = "Normal"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.41'
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '40.87'
$ws.Range("D28").Style = "Normal"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '172.63'
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '21.44'
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0909'
$ws.Range("D33").Style = "Normal"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '4.67'
$ws.Range("D36").Style = "Normal"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0358'
$ws.Range("D38").Style = "Normal"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.69'
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '78.11'
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '14.02'
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.240'
$ws.Range("D43").Style = "Normal"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '8.73'
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '104.44'
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0998'
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.25'
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.428.75'
$ws.Range("D51").Style = "Normal"

# Coin name (B), Link (C) and Volume(1h) (E) column updates
$ws.Range("E2").Value = '  -1.09%  '
$ws.Range("E3").Value = '  -0.15%  '
$ws.Range("E4").Value = '  -0.13%  '
$ws.Range("E5").Value = '  -0.38%  '
$ws.Range("E6").Value = '  -0.90%  '
$ws.Range("E7").Value = '  +0.00%  '
$ws.Range("E8").Value = '  -0.23%  '
$ws.Range("E9").Value = '  -0.56%  '
$ws.Range("E10").Value = '  -0.61%  '
$ws.Range("E11").Value = '  -0.43%  '
$ws.Range("E12").Value = '  +5.84%  '
$ws.Range("E13").Value = '  +1.56%  '
$ws.Range("E14").Value = '  +0.35%  '
$ws.Range("E15").Value = '  -0.36%  '
$ws.Range("E16").Value = '  +2.37%  '
$ws.Range("E17").Value = '  -0.41%  '
$ws.Range("E18").Value = '  -0.67%  '
$ws.Range("E19").Value = '  -0.85%  '
$ws.Range("E20").Value = '  +4.97%  '
$ws.Range("E21").Value = '  -1.04%  '
$ws.Range("E22").Value = '  -0.14%  '
$ws.Range("E23").Value = '  +0.04%  '
$ws.Range("E24").Value = '  +0.82%  '
$ws.Range("E25").Value = '  +1.84%  '
$ws.Range("E26").Value = '  +1.29%  '
$ws.Range("E27").Value = '  +0.25%  '
$ws.Range("E28").Value = '  -6.95%  '
$ws.Range("E29").Value = '  -2.35%  '
$ws.Range("E30").Value = '  -0.92%  '
$ws.Range("E31").Value = '  -2.01%  '
$ws.Range("E32").Value = '  -0.89%  '
$ws.Range("E33").Value = '  -1.85%  '
$ws.Range("E34").Value = '  +5.62%  '
$ws.Range("E35").Value = '  +0.97%  '
$ws.Range("E36").Value = '  -1.06%  '
$ws.Range("E37").Value = '  +1.40%  '
$ws.Range("E38").Value = '  +1.22%  '
$ws.Range("E39").Value = '  -4.14%  '
$ws.Range("E40").Value = '  +13.02%  '
$ws.Range("E41").Value = '  +3.68%  '
$ws.Range("E42").Value = '  +5.99%  '
$ws.Range("E43").Value = '  -0.76%  '
$ws.Range("E44").Value = '  +5.47%  '
$ws.Range("E45").Value = '  -0.11%  '
$ws.Range("E46").Value = '  -1.93%  '
$ws.Range("E47").Value = '  +0.01%  '
$ws.Range("E48").Value = '  +3.70%  '
$ws.Range("B49").Value = 'Cronos'
$ws.Range("C49").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("E49").Value = '  -0.31%  '
$ws.Range("B50").Value = 'TrustWalletToken'
$ws.Range("C50").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("E50").Value = '  +1.96%  '
$ws.Range("B51").Value = 'Maker'
$ws.Range("C51").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("E51").Value = '  +2.30%  '
